# Update symbol list (cryptos) - price refresh + row 42/43 swap
# Commit: "Updated symbol list on Fri Dec 30 21:22:16 UTC 2022 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking value into a cell while keeping it stored
# as TEXT (matches the source data, which is all inlineStr/text). Assigning
# a plain numeric-looking string via .Value makes Excel COM auto-convert the
# cell to a number, so we briefly force Text format, assign, then restore the
# cell's original (General) style by copying it back from an untouched cell
# so no stray style index lingers on the cell.
$blankStyleCell = $ws.Range("D6")

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $blankStyleCell.Style
}

# --- Column D (Price) refresh ---
Set-TextValue "D2"  "245.01"
Set-TextValue "D3"  "25.19"
Set-TextValue "D4"  "5.009"
Set-TextValue "D5"  "0.05613"
Set-TextValue "D7"  "3.006"
Set-TextValue "D8"  "0.8133"
Set-TextValue "D9"  "0.8387"
Set-TextValue "D10" "0.1336"
Set-TextValue "D11" "0.06944"
Set-TextValue "D12" "0.03263"
Set-TextValue "D13" "0.02838"
Set-TextValue "D14" "0.09395"
Set-TextValue "D15" "0.001512"
Set-TextValue "D16" "0.0005961"
Set-TextValue "D17" "0.006096"
Set-TextValue "D18" "3.501"
Set-TextValue "D21" "0.1291"
Set-TextValue "D22" "3.737"
Set-TextValue "D24" "0.1368"
Set-TextValue "D25" "0.001244"
Set-TextValue "D26" "0.004525"
Set-TextValue "D27" "0.00009693"
Set-TextValue "D28" "0.0001937"
Set-TextValue "D40" "0.03662"
Set-TextValue "D41" "0.1365"
Set-TextValue "D44" "0.008172"
Set-TextValue "D45" "0.00005286"
Set-TextValue "D46" "0.00000000749"
Set-TextValue "D47" "0.2256"
Set-TextValue "D49" "0.00002096"
Set-TextValue "D50" "0.0001997"

# --- Row 42 / Row 43 data swap (Coin, Link, Price, Volume columns B:E) ---
# Row 42 becomes CEJI, Row 43 becomes KickToken (values updated slightly too)
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.002720"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D43" "0.006226"
$ws.Range("E43").Value = "42KickTokenKICK"
